$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new %Diff columns, styled like the existing Diff headers (AC1:AE1)
$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

# Copy the header style (bold, border, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AF1:AH1").PasteSpecial(-4122)

# %Diff = Diff / Ytrue * 100, for each row and each outcome (H, D, A)
for ($row = 2; $row -le 4; $row++) {
    $diffH = $ws.Cells.Item($row, 29).Value2   # AC
    $diffD = $ws.Cells.Item($row, 30).Value2   # AD
    $diffA = $ws.Cells.Item($row, 31).Value2   # AE

    $trueH = $ws.Cells.Item($row, 26).Value2   # Z
    $trueD = $ws.Cells.Item($row, 27).Value2   # AA
    $trueA = $ws.Cells.Item($row, 28).Value2   # AB

    $ws.Cells.Item($row, 32).Value = $diffH / $trueH * 100   # AF
    $ws.Cells.Item($row, 33).Value = $diffD / $trueD * 100   # AG
    $ws.Cells.Item($row, 34).Value = $diffA / $trueA * 100   # AH
}
